# Scheduled-runner style update of market-board-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# eight job sheets, reflecting refreshed price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31: Hush Little Wailer
$ws.Range("H31").Value2 = 2500
$ws.Range("I31").Value2 = 2500
$ws.Range("K31").Value2 = 7500
$ws.Range("M31").Value2 = -7270

# Row 32: Automata for the People
$ws.Range("H32").Value2 = 3502.0833
$ws.Range("I32").Value2 = 2950
$ws.Range("J32").Value2 = 3778.125
$ws.Range("K32").Value2 = 2950
$ws.Range("L32").Value2 = 3778.125
$ws.Range("M32").Value2 = -2624
$ws.Range("N32").Value2 = -4430.125

# Row 41: The Write Stuff
$ws.Range("H41").Value2 = 9615698
$ws.Range("I41").Value2 = 13889100
$ws.Range("J41").Value2 = 541.75
$ws.Range("K41").Value2 = 13889100
$ws.Range("L41").Value2 = 541.75
$ws.Range("M41").Value2 = -13888660
$ws.Range("N41").Value2 = -1421.75

# Row 51: A Bile Business
$ws.Range("H51").Value2 = 14255.5
$ws.Range("I51").Value2 = 24282.666
$ws.Range("J51").Value2 = 9958.143
$ws.Range("K51").Value2 = 24282.666
$ws.Range("L51").Value2 = 9958.143
$ws.Range("M51").Value2 = -23798.666
$ws.Range("N51").Value2 = -10926.143

# Row 70: Consecrating Congregation
$ws.Range("H70").Value2 = 72225300
$ws.Range("I70").Value2 = 35716370
$ws.Range("J70").Value2 = 104170610
$ws.Range("K70").Value2 = 107149110
$ws.Range("L70").Value2 = 312511830
$ws.Range("M70").Value2 = -107148840
$ws.Range("N70").Value2 = -312512370

# Row 73: Curbing the Contagion (L)
$ws.Range("H73").Value2 = 72225300
$ws.Range("I73").Value2 = 35716370
$ws.Range("J73").Value2 = 104170610
$ws.Range("K73").Value2 = 107149110
$ws.Range("L73").Value2 = 312511830
$ws.Range("M73").Value2 = -107148174
$ws.Range("N73").Value2 = -312513702

# Row 76: Warding Off Temptation
$ws.Range("H76").Value2 = 10338.3
$ws.Range("I76").Value2 = 11684.875
$ws.Range("K76").Value2 = 11684.875
$ws.Range("M76").Value2 = -11369.875

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value2 = 10338.3
$ws.Range("I79").Value2 = 11684.875
$ws.Range("K79").Value2 = 11684.875
$ws.Range("M79").Value2 = -10592.875

# Row 88: The Grave of Hemlock Groves
$ws.Range("H88").Value2 = 12377484
$ws.Range("I88").Value2 = 37040104
$ws.Range("J88").Value2 = 46174.223
$ws.Range("K88").Value2 = 37040104
$ws.Range("L88").Value2 = 46174.223
$ws.Range("M88").Value2 = -37039698
$ws.Range("N88").Value2 = -46986.223

# Row 91: Dappling the Highlands (L)
$ws.Range("H91").Value2 = 12377484
$ws.Range("I91").Value2 = 37040104
$ws.Range("J91").Value2 = 46174.223
$ws.Range("K91").Value2 = 37040104
$ws.Range("L91").Value2 = 46174.223
$ws.Range("M91").Value2 = -37038700
$ws.Range("N91").Value2 = -48982.223

# Row 107: Another Man's Ink
$ws.Range("H107").Value2 = 15001262
$ws.Range("I107").Value2 = 5953297
$ws.Range("J107").Value2 = 62503076
$ws.Range("K107").Value2 = 5953297
$ws.Range("L107").Value2 = 62503076
$ws.Range("M107").Value2 = -5951377
$ws.Range("N107").Value2 = -62506916

# Row 116: Growing Up
$ws.Range("H116").Value2 = 27786940
$ws.Range("I116").Value2 = 62506124
$ws.Range("J116").Value2 = 11595
$ws.Range("K116").Value2 = 62506124
$ws.Range("L116").Value2 = 11595
$ws.Range("M116").Value2 = -62502682
$ws.Range("N116").Value2 = -18479

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value2 = 3090.4348
$ws.Range("I137").Value2 = 3536.875
$ws.Range("K137").Value2 = 10610.625
$ws.Range("M137").Value2 = -8060.625

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value2 = 3368.423
$ws.Range("I2").Value2 = 2212.4666
$ws.Range("K2").Value2 = 2212.4666
$ws.Range("M2").Value2 = -2099.4666

# Row 32: Ingot We Trust
$ws.Range("H32").Value2 = 2409520.2
$ws.Range("I32").Value2 = 2504750.8
$ws.Range("K32").Value2 = 2504750.8
$ws.Range("M32").Value2 = -2504463.8

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value2 = 12346500
$ws.Range("I110").Value2 = 677.4091
$ws.Range("K110").Value2 = 677.4091
$ws.Range("M110").Value2 = 1367.5909

# Row 116: No Scope
$ws.Range("H116").Value2 = 3368.423
$ws.Range("I116").Value2 = 2212.4666
$ws.Range("K116").Value2 = 2212.4666
$ws.Range("M116").Value2 = 81.5333999999998

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value2 = 3565.4
$ws.Range("I132").Value2 = 1425.6389
$ws.Range("K132").Value2 = 4276.9167
$ws.Range("M132").Value2 = -1746.9167

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value2 = 3368.423
$ws.Range("I3").Value2 = 2212.4666
$ws.Range("K3").Value2 = 2212.4666
$ws.Range("M3").Value2 = -2098.4666

# Row 20: Smelt and Dealt
$ws.Range("H20").Value2 = 9806214
$ws.Range("J20").Value2 = 3218.4
$ws.Range("L20").Value2 = 3218.4
$ws.Range("N20").Value2 = -3712.4

# Row 86: Through Thick and Thin
$ws.Range("H86").Value2 = 34725244
$ws.Range("I86").Value2 = 14708434
$ws.Range("K86").Value2 = 14708434
$ws.Range("M86").Value2 = -14707311

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value2 = 34725244
$ws.Range("I89").Value2 = 14708434
$ws.Range("K89").Value2 = 73542170
$ws.Range("M89").Value2 = -73536554

# Row 94: High Steal
$ws.Range("H94").Value2 = 1708.875
$ws.Range("I94").Value2 = 749.65
$ws.Range("J94").Value2 = 6505
$ws.Range("K94").Value2 = 749.65
$ws.Range("L94").Value2 = 6505
$ws.Range("M94").Value2 = -298.65
$ws.Range("N94").Value2 = -7407

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value2 = 6565.3335
$ws.Range("I16").Value2 = 4400
$ws.Range("K16").Value2 = 4400
$ws.Range("M16").Value2 = -4113

# Row 22: Driving Up the Wall
$ws.Range("H22").Value2 = 349.5
$ws.Range("I22").Value2 = 399.5
$ws.Range("J22").Value2 = 299.5
$ws.Range("K22").Value2 = 399.5
$ws.Range("L22").Value2 = 299.5
$ws.Range("M22").Value2 = -49.5
$ws.Range("N22").Value2 = -999.5

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value2 = 17865394
$ws.Range("I58").Value2 = 71429270
$ws.Range("J58").Value2 = 10765.81
$ws.Range("K58").Value2 = 71429270
$ws.Range("L58").Value2 = 10765.81
$ws.Range("M58").Value2 = -71429067
$ws.Range("N58").Value2 = -11171.81

# Row 107: Built to Last
$ws.Range("H107").Value2 = 1957.5
$ws.Range("I107").Value2 = 1831.75
$ws.Range("K107").Value2 = 1831.75
$ws.Range("M107").Value2 = 88.25

# Row 113: Patient Patients
$ws.Range("H113").Value2 = 6565.3335
$ws.Range("I113").Value2 = 4400
$ws.Range("K113").Value2 = 4400
$ws.Range("M113").Value2 = -2230

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value2 = 4556
$ws.Range("J122").Value2 = 5608.4
$ws.Range("L122").Value2 = 16825.2
$ws.Range("N122").Value2 = -21725.2

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value2 = 5754.2646
$ws.Range("I132").Value2 = 1789.1765
$ws.Range("K132").Value2 = 5367.529500000001
$ws.Range("M132").Value2 = -2837.529500000001

# Row 136: Turali Quality
$ws.Range("H136").Value2 = 17865394
$ws.Range("I136").Value2 = 71429270
$ws.Range("J136").Value2 = 10765.81
$ws.Range("K136").Value2 = 214287810
$ws.Range("L136").Value2 = 32297.43
$ws.Range("M136").Value2 = -214285260
$ws.Range("N136").Value2 = -37397.43

$ws = $wb.Worksheets.Item("CUL")
# Row 3: Trout Fishing in Limsa
$ws.Range("H3").Value2 = 3717
$ws.Range("I3").Value2 = 3717
$ws.Range("K3").Value2 = 11151
$ws.Range("M3").Value2 = -11039

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value2 = 4017.6428
$ws.Range("I80").Value2 = 3216.1428
$ws.Range("K80").Value2 = 3216.1428
$ws.Range("M80").Value2 = -2218.1428

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value2 = 4017.6428
$ws.Range("I83").Value2 = 3216.1428
$ws.Range("K83").Value2 = 16080.714
$ws.Range("M83").Value2 = -11088.714

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value2 = 7596.517
$ws.Range("I113").Value2 = 3749.875
$ws.Range("K113").Value2 = 3749.875
$ws.Range("M113").Value2 = -1579.875

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value2 = 2015810.2
$ws.Range("I122").Value2 = 3152699
$ws.Range("K122").Value2 = 9458097
$ws.Range("M122").Value2 = -9455647

# Row 132: On Board for Lar
$ws.Range("H132").Value2 = 4017.5557
$ws.Range("I132").Value2 = 1562.9565
$ws.Range("J132").Value2 = 8360.308000000001
$ws.Range("K132").Value2 = 4688.8695
$ws.Range("L132").Value2 = 25080.924
$ws.Range("M132").Value2 = -2158.8695
$ws.Range("N132").Value2 = -30140.924

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value2 = 5712.5625
$ws.Range("J7").Value2 = 6000.0835
$ws.Range("L7").Value2 = 6000.0835
$ws.Range("N7").Value2 = -6224.0835

# Row 40: Best Served Toad
$ws.Range("H40").Value2 = 7520.2
$ws.Range("I40").Value2 = 5749.75
$ws.Range("J40").Value2 = 8164
$ws.Range("K40").Value2 = 5749.75
$ws.Range("L40").Value2 = 8164
$ws.Range("M40").Value2 = -5613.75
$ws.Range("N40").Value2 = -8436

# Row 46: Supply Side Logic
$ws.Range("H46").Value2 = 1265.2273
$ws.Range("J46").Value2 = 1433.421
$ws.Range("L46").Value2 = 1433.421
$ws.Range("N46").Value2 = -1809.421

# Row 61: Spelling Me Softly
$ws.Range("H61").Value2 = 4046.1724
$ws.Range("I61").Value2 = 1363.8334
$ws.Range("K61").Value2 = 1363.8334
$ws.Range("M61").Value2 = -1161.8334

# Row 113: Peace in Rest
$ws.Range("H113").Value2 = 4046.1724
$ws.Range("I113").Value2 = 1363.8334
$ws.Range("K113").Value2 = 1363.8334
$ws.Range("M113").Value2 = 806.1666

# Row 126: Battered Books
$ws.Range("H126").Value2 = 5712.5625
$ws.Range("J126").Value2 = 6000.0835
$ws.Range("L126").Value2 = 18000.2505
$ws.Range("N126").Value2 = -22940.2505

# Row 132: Tenets of Tanning
$ws.Range("H132").Value2 = 14714107
$ws.Range("J132").Value2 = 10241.417
$ws.Range("L132").Value2 = 30724.251
$ws.Range("N132").Value2 = -35784.251

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax
$ws.Range("H107").Value2 = 17544700
$ws.Range("I107").Value2 = 579.7778
$ws.Range("J107").Value2 = 33334408
$ws.Range("K107").Value2 = 1739.3334
$ws.Range("L107").Value2 = 100003224
$ws.Range("M107").Value2 = 180.6666
$ws.Range("N107").Value2 = -100007064

# Row 113: A Tender Table
$ws.Range("H113").Value2 = 2388.2856
$ws.Range("I113").Value2 = 1626
$ws.Range("J113").Value2 = 2960
$ws.Range("K113").Value2 = 4878
$ws.Range("L113").Value2 = 8880
$ws.Range("M113").Value2 = -2708
$ws.Range("N113").Value2 = -13220

# Row 126: A Polished Purchase
$ws.Range("H126").Value2 = 3102.0454
$ws.Range("I126").Value2 = 1380.2307
$ws.Range("K126").Value2 = 4140.6921
$ws.Range("M126").Value2 = -1670.6921
